$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.282.11"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.88%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.322.32"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "588.15"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.59%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "183.34"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.60%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.648"
$c.ClearFormats()
$ws.Range("E7").Value = "  +8.12%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("E10").Value = "  +2.11%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.401"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.43%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "3.897.05"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("E13").Value = "  -3.76%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "66.300.67"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.93%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "26.28"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.37%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000163"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.67%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.283.40"
$c.ClearFormats()
$ws.Range("E17").Value = "  -1.84%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "428.78"
$c.ClearFormats()
$ws.Range("E18").Value = "  -1.95%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.54"
$c.ClearFormats()
$ws.Range("E19").Value = "  -2.68%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.21"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.39%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.42"
$c.ClearFormats()
$ws.Range("E21").Value = "  -2.97%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "71.98"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.58%  "

$ws.Range("E23").Value = "  -0.05%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.39%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.461.51"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.97%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.518"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.54%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.206"
$c.ClearFormats()
$ws.Range("E27").Value = "  +7.87%  "

$ws.Range("E28").Value = "  -4.14%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.95"
$c.ClearFormats()
$ws.Range("E29").Value = "  -1.74%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  -0.92%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "22.39"
$c.ClearFormats()
$ws.Range("E32").Value = "  -2.43%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.19"
$c.ClearFormats()
$ws.Range("E34").Value = "  -2.47%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.59"
$c.ClearFormats()
$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("E36").Value = "  -4.17%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "159.72"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("E38").Value = "  -3.89%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.881.87"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("E40").Value = "  -2.40%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.61"
$c.ClearFormats()
$ws.Range("E41").Value = "  -4.60%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.766"
$c.ClearFormats()
$ws.Range("E42").Value = "  -3.46%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.33"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.49%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "40.15"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.28%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0666"
$c.ClearFormats()
$ws.Range("E45").Value = "  -1.36%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.87"
$c.ClearFormats()
$ws.Range("E46").Value = "  -5.93%  "

$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.25"
$c.ClearFormats()
$ws.Range("E48").Value = "  -5.76%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "315.86"
$c.ClearFormats()
$ws.Range("E49").Value = "  -2.75%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0272"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.40%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.ClearFormats()
$ws.Range("E51").Value = "  +5.03%  "
